$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells so numeric-looking strings (e.g. "28.10",
# "0.0000187", dotted thousand-separator prices) keep their exact text
# representation instead of being auto-converted to floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.821.83"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.605.38"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.25"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.34"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.603.23"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.26"
$ws.Range("E11").Value = "  +3.97%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.219.85"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.10"
$ws.Range("E14").Value = "  +2.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.605.40"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.922.29"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.12"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.71"
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.93"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "398.02"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.751.48"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.60"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.20"
$ws.Range("E28").Value = "  +4.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.68"
$ws.Range("E29").Value = "  +28.86%  "
$ws.Range("E30").Value = "  +4.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.68"
$ws.Range("E31").Value = "  +4.22%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.608.38"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.68"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.148"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("E37").Value = "  +7.95%  "
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.11"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "171.95"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0837"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.05"
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.41"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.25"
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.55"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("E49").Value = "  +4.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.455.43"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  +3.30%  "
